# Update French lakes' (geneva, bourget, annecy) spawning.depth.m values
# from 2 to 4 on rows 3-5 (column G).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bio-parameters")

$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 4

# Update the active selection to match the saved view state.
$ws.Range("H14").Select()
